$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 4 down to the new row 5 first, so the new cells
# reuse the existing style records (s="1"/"2"/"3") instead of Excel minting
# brand-new ones for the freshly entered values.
$ws.Range("A4:I4").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122) # xlPasteFormats

# Append new row for LeetCode problem 3443: Maximum Manhattan Distance After K Changes
$ws.Range("A5").Value = 3443
$ws.Range("B5").Value = "Maximum Manhattan Distance After K Changes"
$ws.Range("C5").Value = "math, string, counting"
$ws.Range("D5").Value = "medium"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 45828
$ws.Range("I5").Value = 45828

# Match the row height used by row 4 (same visual size as the other wrapped rows)
$ws.Rows.Item(5).RowHeight = 51

# The sheet's outlineLevelRow high-water mark bumps 3 -> 4 alongside this edit.
$ws.Rows.Item(7).OutlineLevel = 4
$ws.Rows.Item(7).EntireRow.Delete()

# Leave the cursor on D4, matching the post-edit selection state.
$ws.Range("D4").Select()
